$wb = $excel.ActiveWorkbook

# --- Step 1: swap the two sheet tab names -----------------------------
# Before: Worksheets(1) = "hotel_info", Worksheets(2) = "review_info"
# After:  Worksheets(1) = "review_info", Worksheets(2) = "hotel_info"
$wsFirst  = $wb.Worksheets.Item(1)
$wsSecond = $wb.Worksheets.Item(2)

$wsFirst.Name  = "__swap_tmp__"
$wsSecond.Name = "hotel_info"
$wsFirst.Name  = "review_info"

# --- Step 2: rebuild the "review_info" sheet (now Worksheets(1)) ------
# It becomes a header-only sheet (25 columns), the old hotel_info data
# that used to live here is gone.
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Cells.Clear()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $wsReview.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- Step 3: rebuild the "hotel_info" sheet (now Worksheets(2)) -------
# It gets the hotel header row plus the single data row, with a new
# "State" column inserted right after "Hotel_Name".
$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsHotel.Cells.Clear()

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $wsHotel.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$wsHotel.Cells.Item(2, 1).Value = 55074
$wsHotel.Cells.Item(2, 2).Value = "Homewood Suites Covington"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"
$wsHotel.Cells.Item(2, 4).Value = "Covington"
$wsHotel.Cells.Item(2, 5).Value = 70433
$wsHotel.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g40095-d618912-Reviews-Homewood_Suites_Covington-Covington_Louisiana.html"
$wsHotel.Cells.Item(2, 7).Value = "Homewood Suites Covington"

# These three look numeric but are stored as text in the source data, so
# force a text format before assigning (otherwise Excel auto-detects them
# as numbers).
$wsHotel.Range("H2:J2").NumberFormat = "@"
$wsHotel.Cells.Item(2, 8).Value = "345"
$wsHotel.Cells.Item(2, 9).Value = "5"
$wsHotel.Cells.Item(2, 10).Value = "347"
